# Unit Three.pptx - "new updates for chapter three"
#
# Slide 46, shape "Content Placeholder 2": the run containing
# "GitHub link " becomes "GitHub Link" and is turned into a hyperlink
# (<a:hlinkClick r:id="rId3"/> added to that run's rPr).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(46)

# Find the "Content Placeholder 2" shape (body placeholder) on the slide.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tr = $shape.TextFrame.TextRange

# Locate the "GitHub link " run inside the full text of the placeholder
# instead of hard-coding character offsets.
$needle = "GitHub link "
$fullText = $tr.Text
$pos = $fullText.IndexOf($needle)

$start = $pos + 1                   # PowerPoint TextRange is 1-based
$len = $needle.Length
$run = $tr.Characters($start, $len)

# Update the visible text and turn it into a hyperlink.
$run.Text = "GitHub Link"
$run.ActionSettings(1).Hyperlink.Address = "https://github.com/"
